$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.437.05"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "1.553.70"
$ws.Range("E3").Value = "  -1.90%  "

$ws.Range("D5").Value = "210.47"
$ws.Range("E5").Value = "  -1.48%  "

$ws.Range("E6").Value = "  -1.66%  "

$ws.Range("E7").Value = "  -0.29%  "

$ws.Range("D8").Value = "24.14"
$ws.Range("E8").Value = "  +0.91%  "

$ws.Range("E9").Value = "  -1.81%  "

$ws.Range("E10").Value = "  -1.32%  "

$ws.Range("E11").Value = "  -0.46%  "

$ws.Range("D12").Value = "1.773.31"
$ws.Range("E12").Value = "  -2.05%  "

$ws.Range("D13").Value = "1.544.15"
$ws.Range("E13").Value = "  -2.67%  "

$ws.Range("D14").Value = "28.444.81"
$ws.Range("E14").Value = "  -0.06%  "

$ws.Range("E15").Value = "  -2.10%  "

$ws.Range("E16").Value = "  -1.80%  "

$ws.Range("D17").Value = "61.06"
$ws.Range("E17").Value = "  -1.76%  "

$ws.Range("D18").Value = "228.88"
$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").Value = "7.32"
$ws.Range("E19").Value = "  -1.55%  "

$ws.Range("D20").Value = "0.0₃0673"
$ws.Range("E20").Value = "  -2.43%  "

$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("E22").Value = "  -0.49%  "

$ws.Range("D23").Value = "8.91"
$ws.Range("E23").Value = "  -2.55%  "

$ws.Range("D24").Value = "2.02"
$ws.Range("E24").Value = "  -2.08%  "

$ws.Range("D25").Value = "151.18"
$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").Value = "14.75"
$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("E29").Value = "  -3.25%  "

$ws.Range("E30").Value = "  -3.25%  "

$ws.Range("E31").Value = "  -4.57%  "

$ws.Range("E32").Value = "  -1.87%  "

$ws.Range("D33").Value = "1.382.69"
$ws.Range("E33").Value = "  -0.95%  "

$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("D35").Value = "1.05"
$ws.Range("E35").Value = "  -2.63%  "

$ws.Range("D36").Value = "1.48"
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("D37").Value = "2.29"
$ws.Range("E37").Value = "  -2.99%  "

$ws.Range("E38").Value = "  -3.00%  "

$ws.Range("E39").Value = "  -2.52%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "1.92"
$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "0.510"
$ws.Range("E41").Value = "  -2.30%  "

$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("D43").Value = "0.773"
$ws.Range("E43").Value = "  -2.15%  "

$ws.Range("D44").Value = "0.0459"
$ws.Range("E44").Value = "  -0.37%  "

$ws.Range("D45").Value = "5.36"
$ws.Range("E45").Value = "  -1.55%  "

$ws.Range("D46").Value = "61.76"

$ws.Range("D47").Value = "1.687.45"
$ws.Range("E47").Value = "  -1.99%  "

$ws.Range("D48").Value = "0.876"
$ws.Range("E48").Value = "  -8.85%  "

$ws.Range("E49").Value = "  -1.60%  "

$ws.Range("D50").Value = "43.12"
$ws.Range("E50").Value = "  +8.76%  "

$ws.Range("D51").Value = "0.0₆0101"
$ws.Range("E51").Value = "  -1.67%  "
